$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $r = $ws.Range($Range)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

Set-TextValue "D2" "26.178.05"
Set-TextValue "E2" "  +0.45%  "
Set-TextValue "D3" "1.654.49"
Set-TextValue "E3" "  +0.08%  "
Set-TextValue "E4" "  +0.43%  "
Set-TextValue "D5" "217.73"
Set-TextValue "E5" "  -0.32%  "
Set-TextValue "D6" "0.5315"
Set-TextValue "E6" "  +0.57%  "
Set-TextValue "D8" "0.2627"
Set-TextValue "E8" "  +0.20%  "
Set-TextValue "E9" "  +0.20%  "
Set-TextValue "D10" "20.39"
Set-TextValue "E10" "  -0.24%  "
Set-TextValue "D11" "0.07813"
Set-TextValue "E11" "  +0.93%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.701.65"
Set-TextValue "E12" "  +1.43%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.518"
Set-TextValue "E13" "  +0.65%  "
Set-TextValue "D14" "1.881.28"
Set-TextValue "E14" "  +0.19%  "
Set-TextValue "E15" "  +0.27%  "
Set-TextValue "D16" "0.0₅8151"
Set-TextValue "E16" "  +0.55%  "
Set-TextValue "D17" "65.41"
Set-TextValue "E17" "  +0.40%  "
Set-TextValue "D18" "26.137.18"
Set-TextValue "E18" "  +0.21%  "
Set-TextValue "E19" "  +0.42%  "
Set-TextValue "E20" "  +0.35%  "
Set-TextValue "D21" "191.06"
Set-TextValue "E21" "  -1.30%  "
Set-TextValue "E22" "  +0.14%  "
Set-TextValue "D23" "6.012"
Set-TextValue "E23" "  +0.18%  "
Set-TextValue "D24" "1.008"
Set-TextValue "E24" "  +0.47%  "
Set-TextValue "D25" "145.30"
Set-TextValue "E25" "  +3.95%  "
Set-TextValue "D26" "0.1218"
Set-TextValue "E26" "  -2.12%  "
Set-TextValue "D27" "7.189"
Set-TextValue "E27" "  -1.21%  "
Set-TextValue "E28" "  -2.08%  "
Set-TextValue "D29" "1.475"
Set-TextValue "E29" "  +4.47%  "
Set-TextValue "D30" "0.05748"
Set-TextValue "E30" "  -3.40%  "
Set-TextValue "D31" "1.273"
Set-TextValue "E31" "  -0.20%  "
Set-TextValue "D32" "3.550"
Set-TextValue "E32" "  +1.39%  "
Set-TextValue "D33" "3.263"
Set-TextValue "E33" "  +0.39%  "
Set-TextValue "D34" "1.586"
Set-TextValue "E34" "  +2.81%  "
Set-TextValue "E35" "  +1.91%  "
Set-TextValue "E36" "  +0.32%  "
Set-TextValue "D37" "0.9485"
Set-TextValue "E37" "  +0.52%  "
Set-TextValue "D38" "0.5747"
Set-TextValue "E38" "  +1.72%  "
Set-TextValue "E39" "  -0.71%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D40" "0.8498"
Set-TextValue "E40" "  +0.29%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "5.789"
Set-TextValue "E41" "  -1.70%  "
Set-TextValue "E42" "  +0.46%  "
Set-TextValue "D43" "1.040.69"
Set-TextValue "E43" "  +3.29%  "
Set-TextValue "D44" "103.76"
Set-TextValue "E44" "  +2.81%  "
Set-TextValue "D45" "1.793.91"
Set-TextValue "E45" "  +0.02%  "
Set-TextValue "D46" "56.72"
Set-TextValue "E46" "  -0.43%  "
Set-TextValue "E47" "  -1.66%  "
Set-TextValue "D48" "1.002"
Set-TextValue "E48" "  -0.26%  "
Set-TextValue "D49" "0.4358"
Set-TextValue "E49" "  +1.65%  "
Set-TextValue "D50" "7.865"
Set-TextValue "E50" "  +0.19%  "
Set-TextValue "D51" "0.05153"
Set-TextValue "E51" "  +0.07%  "
